$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated values for rows 2-9 (columns B:F), and rank values in column G
$data = @(
    @{ row = 2;  B = -0.3989960852331779; C = 1.622347057549135;  D = 3.410040606493075;   E = 1.846629526053636;  F = 1.871071232177386;  G = 14 }
    @{ row = 3;  B = -0.07368367304384635; C = 2.01870334081424;  D = 6.614665517794586;   E = 2.571899204439121;  F = 2.675818738311348;  G = 13 }
    @{ row = 4;  B = -0.3412105586050318; C = 1.726870111009629;  D = 4.898523348216152;   E = 2.213260795346123;  F = 2.284039174058081;  G = 12 }
    @{ row = 5;  B = -0.4693883225059389; C = 2.197338286322757;  D = 7.681313688700947;   E = 2.771518300264486;  F = 2.864801410313742;  G = 11 }
    @{ row = 6;  B = -0.7369237006858437; C = 1.317817770111458;  D = 2.080685236058409;   E = 1.442458053483154;  F = 1.30708696277294;   G = 10 }
    @{ row = 7;  B = -0.2355834118156171; C = 1.992865734783982;  D = 4.408914335405655;   E = 2.099741492518937;  F = 2.213050302279264;  G = 9 }
    @{ row = 8;  B = 0.04677381843788716; C = 1.625585584899291;  D = 3.112956264121649;   E = 1.764357181559802;  F = 1.932077164306971;  G = 6 }
    @{ row = 9;  B = 0.8102998313761515;  C = 0.9132997309681818; D = 1.270800624281741;   E = 1.127297930576359;  F = 0.9598553074970639; G = 3 }
)

foreach ($d in $data) {
    $r = $d.row
    $ws.Cells.Item($r, 2).Value = $d.B
    $ws.Cells.Item($r, 3).Value = $d.C
    $ws.Cells.Item($r, 4).Value = $d.D
    $ws.Cells.Item($r, 5).Value = $d.E
    $ws.Cells.Item($r, 6).Value = $d.F
    $ws.Cells.Item($r, 7).Value = $d.G
}

# New row 10 - add the Q8 label and its values
$ws.Cells.Item(10, 1).Value = "Q8"
$ws.Cells.Item(9, 1).Copy()
$ws.Cells.Item(10, 1).PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Cells.Item(10, 2).Value = 0.1034734828819666
$ws.Cells.Item(10, 3).Value = 0.1034734828819666
$ws.Cells.Item(10, 4).Value = 0.01070676165972463
$ws.Cells.Item(10, 5).Value = 0.1034734828819666
$ws.Cells.Item(10, 7).Value = 1
